# Completanto funcionalidad de pagos pendientes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (week 33): MIERCOLES (D) now paid, update TOTAL
$ws.Range("D2").Value = "SI"
$ws.Range("G2").Value = 240

# Row 3 (week 34): JUEVES (E) and VIERNES (F) now paid, update TOTAL
$ws.Range("E3").Value = "SI"
$ws.Range("F3").Value = "SI"
$ws.Range("G3").Value = 240

# Row 4 (week 35): LUNES (B) and MIERCOLES (D) now paid, update TOTAL
$ws.Range("B4").Value = "SI"
$ws.Range("D4").Value = "SI"
$ws.Range("G4").Value = 160
